$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated odds values per the diff (cell -> new value)
$updates = @{
    "M3" = 1.07
    "O3" = 1.41
    "P3" = 2.7
    "M4" = 1.05
    "O4" = 1.41
    "P4" = 2.7
    "G5" = 2.55
    "H5" = 2.88
    "G6" = 1.8
    "R6" = 1.57
    "G7" = 2.6
    "I7" = 2.6
    "M8" = 1.03
    "O8" = 1.25
    "G9" = 3.7
    "H9" = 3.6
    "S9" = 1.33
    "T9" = 3.25
    "W9" = 13
    "AC9" = 13
    "AE9" = 15
    "AG9" = 151
    "AK9" = 15
    "AT9" = 3.25
    "G11" = 1.66
    "M11" = 1.02
    "O11" = 1.15
    "G13" = 4.1
    "H13" = 4.33
    "I13" = 1.7
    "J13" = 4
    "L13" = 2.2
    "U13" = 1.44
    "V13" = 2.63
    "X13" = 26
    "Y13" = 15
    "AA13" = 26
    "AB13" = 26
    "AD13" = 9.5
    "AF13" = 34
    "AI13" = 12
    "AK13" = 15
    "AL13" = 12
    "AO13" = 21
    "AV13" = 41
    "AW13" = 4.33
    "AX13" = 8.5
    "AZ13" = 23
    "BA13" = 34
    "M14" = 1.02
    "N14" = 21
    "K15" = 2.37
    "M17" = 1.01
    "O17" = 1.1
    "Q17" = 1.41
    "M18" = 1.03
    "O18" = 1.25
    "M19" = 1.02
    "N19" = 15
    "O19" = 1.19
    "M20" = 1.01
    "O20" = 1.11
    "O21" = 1.07
    "M22" = 1.03
    "O22" = 1.19
    "P22" = 4
    "R23" = 1.57
    "R24" = 1.6
    "J26" = 2.88
    "Q26" = 1.8
    "G27" = 2.05
    "H27" = 3
    "I27" = 3.9
    "M27" = 1.1
    "N27" = 7
    "Q27" = 2.35
    "W27" = 6.5
    "X27" = 9
    "Z27" = 19
    "AC27" = 7
    "AF27" = 51
    "AL27" = 34
    "AN27" = 4
    "AO27" = 12
    "AX27" = 21
    "BA27" = 101
    "Q28" = 1.92
    "R28" = 1.82
    "G29" = 3.6
    "I29" = 1.91
    "J29" = 4.05
    "K29" = 2.15
    "L29" = 2.5
    "N29" = 7.8
    "P29" = 3.5
    "Q29" = 1.78
    "T29" = 2.8
    "W29" = 11.25
    "X29" = 20
    "Y29" = 12
    "Z29" = 50
    "AA29" = 32
    "AC29" = 7.8
    "AD29" = 6.9
    "AF29" = 55
    "AI29" = 9.75
    "AJ29" = 8.25
    "AK29" = 17
    "AL29" = 14.5
    "AN29" = 5.5
    "AO29" = 20
    "AP29" = 26
    "AQ29" = 110
    "AR29" = 150
    "AS29" = 350
    "AT29" = 2.8
    "AW29" = 3.85
    "G32" = 2.3
    "I32" = 2.8
    "L32" = 3.25
    "Q32" = 1.54
    "W32" = 11
    "AE32" = 12
    "AL32" = 21
    "AP32" = 19
    "BC32" = 351
    "R33" = 1.5
    "U35" = 1.91
    "V35" = 1.8
    "U36" = 1.73
    "V38" = 1.73
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
